# Generate Report for Handoff
# The 7827cd75-5671-4ec7-a4ad-afa48cad22e8.md file has been newly handed off
# (again) since its last handback is stale vs. the latest source revision.
# This updates the "Ready for handoff" status + timestamps on the Overview,
# zh-cn and de-de sheets, and records the staleness warning in "Error Detail".

$wb = $excel.ActiveWorkbook

$readyStatus = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/06c895eb5f1db3805ac7e2fa68b5549d219c3155/e2e/7827cd75-5671-4ec7-a4ad-afa48cad22e8.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b982ee8b664dc305895eacb16ac62d012c397240/e2e/7827cd75-5671-4ec7-a4ad-afa48cad22e8.md."

# ---------------------------------------------------------------------------
# Overview sheet: row 3 is the 7827cd75-....md entry. Mark both language
# columns (zh-cn / de-de) "Ready for handoff" and bump the generate date.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $readyStatus
$wsOverview.Range("F3").Value = $readyStatus
$wsOverview.Range("G3").Value = "2016-09-09 12:24:44"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Row 2 (6a6e98f5-....md) simply flips to "Ready for handoff".
$wsZhCn.Range("C2").Value = $readyStatus
# Row 3 (7827cd75-....md) flips too, gets a new Latest Handoff Datetime and
# the staleness warning in Error Detail.
$wsZhCn.Range("C3").Value = $readyStatus
$wsZhCn.Range("H3").Value = "2016-09-09 12:24:33"
$wsZhCn.Range("P3").Value = $errorDetail
# Widen the Error Detail column so the long message is readable.
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $readyStatus
$wsDeDe.Range("C3").Value = $readyStatus
$wsDeDe.Range("H3").Value = "2016-09-09 12:24:44"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
